# Glycans_Defaults.xlsx parameter-table rewrite
# ------------------------------------------------
# The sheet lists IsoMatchMS default parameters (Parameter | Default | Update |
# Description). MatchingAlgorithm + MinimumAbundance are removed, the rest of the
# table shifts up, and three new parameters are introduced: AbundanceThreshold,
# IsotopeMinimum (replacing IsotopeRange) and IsotopingAlgorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 1-4 (header, MZRange, NoiseFilter, Charges) are unchanged.

$ws.Range("A2").Value = 'MZRange'
$ws.Range("B2").Value = '0-6000'
$ws.Range("C2").Value = 'Everytime'
$ws.Range("D2").Value = 'A range of MZ values to filter the data by. It is highly recommended that users visualize the spectra first to determine a reasonable cutoff range.'

$ws.Range("A3").Value = 'NoiseFilter'
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 'Everytime'
$ws.Range("D3").Value = 'An abundance (every peak is scaled to the largest peak) cutoff for peaks. A reasonable value should be in the 2.5 - 5.0% range. Default is 5%.'

$ws.Range("A4").Value = 'Charges'
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 'Everytime'
$ws.Range("D4").Value = 'The range of charges to test. List charges separated by a comma'

$ws.Range("A5").Value = 'AbundanceThreshold'
$ws.Range("B5").Value = 50
$ws.Range("C5").Value = 'Occasionally'
$ws.Range("D5").Value = 'The +/- percent abundance an isotope peak can vary and still be considered a match. If 50%, and the calculated abundance is 3, the matched abundance can vary from 1.5-4'

$ws.Range("A6").Value = 'CorrelationMinimum'
$ws.Range("B6").Value = 0.7
$ws.Range("C6").Value = 'Occasionally'
$ws.Range("D6").Value = 'The minimum correlation value to consider when generating the trelliscope display'

$ws.Range("A7").Value = 'PPMThreshold'
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = 'Occasionally'
$ws.Range("D7").Value = 'The maximum m/z error permitted. '

$ws.Range("A8").Value = 'AdductLabels'
$ws.Range("B8").Value = 'proton,sodium'
$ws.Range("C8").Value = 'Occasionally'
$ws.Range("D8").Value = 'Labels for the Adduct Masses. Should be separated by a comma with no space (ex. proton,sodium)'

$ws.Range("A9").Value = 'AdductMasses'
$ws.Range("B9").Value = '1.00727647,22.98977'
$ws.Range("C9").Value = 'Occasionally'
$ws.Range("D9").Value = 'Masses for the Adducts. Should be separated by a comma with no space (ex. 1.00727647,22.98977)'

$ws.Range("A10").Value = 'AddMAI'
$ws.Range("B10").Value = $false
$ws.Range("C10").Value = 'Occasionally'
$ws.Range("D10").Value = 'Add most abundant isotope to the molecular formula calculation step. Warning: This will slow down the tool. '

$ws.Range("A11").Value = 'IsotopeMinimum'
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 'Rarely'
$ws.Range("D11").Value = 'The minimum number of isotopes to consider. We recommend 5 for intact proteomics, and 2 or 3 otherwise. '

$ws.Range("A12").Value = 'PlottingWindow'
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 'Rarely'
$ws.Range("D12").Value = 'The -/+ m/z value on either side of the matched spectra plot. Default is 2 m/z.'

$ws.Range("A13").Value = 'IsotopingAlgorithm'
$ws.Range("B13").Value = 'Rdisop'
$ws.Range("C13").Value = 'Rarely'
$ws.Range("D13").Value = 'Either "Rdisop" or "isopat". "Rdisop" is more accurate and recommended, though may crash on windows OS. "isopat" may then be used as an alternative. '

# B11 (IsotopeMinimum default) reverts to a plain General number format
# instead of the old date-ish format that used to live on that style slot.
$ws.Range("B11").NumberFormat = "General"

# These description cells were typed fresh (no inherited row styling),
# so clear their formatting like the source workbook does.
$ws.Range("D5").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("A13:D13").ClearFormats()

# New trailing blank row 14, styled like the rest of the table but empty.
$ws.Range("A1:D1").Copy()
$ws.Range("A14:D14").PasteSpecial(-4122)
$ws.Range("A14:D14").ClearContents()
$excel.CutCopyMode = $false

# Restore the active selection to where the editor left off.
$ws.Range("B11").Select()

